$d = $word.ActiveDocument

# Helper: append a new text segment immediately after $prevRange, forcing it
# into its own <w:r> run (rather than being silently re-coalesced into the
# preceding run) while still carrying the Times New Roman run formatting
# used throughout this document.
function New-Segment($prevRange, $text) {
    $r = $word.ActiveDocument.Range($prevRange.End, $prevRange.End)
    $r.InsertAfter($text)
    $r.Font.Name = "Times New Roman"
    $r.Font.NameBi = "Times New Roman"
    $r.Font.Bold = 1
    $r.Font.Bold = 0
    return $r
}

# --- Paragraph 1: "<<Communication Issues>>" body text ---------------------
$p1 = $d.Content
$p1.Find.Execute( `
    "We were able to communicate with all group members, however, Sebastian is unable to attend practical sessions because he lives in London. He is reachable via our MS Teams group chat and we have had group meetings through that.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "We were able to communicate with all group members", 2)

$p1b = New-Segment $p1   " (and we distributed the tasks between us)"
$p1c = New-Segment $p1b  ", "
$p1d = New-Segment $p1c  "except"
$p1e = New-Segment $p1d  " Sebastian "
$p1f = New-Segment $p1e  "who does not reply"
$p1g = New-Segment $p1f  " to "
$p1h = New-Segment $p1g  "emails or messages on MS Teams. "

# --- Paragraph 2: "<<Reflections>>" body text -------------------------------
$p2 = $d.Content
$p2.Find.Execute( `
    "We had no issues with setting up the working environment - we Atlassian/Jira for the board software, and other than that everything seems to be going ok.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "We had no issues with setting up the working environment - we ", 2)

$p2b = New-Segment $p2   "use "
$p2c = New-Segment $p2b  "Atlassian/Jira for the board software and other than that everything seems to be going ok"
$p2d = New-Segment $p2c  " so far."

Write-Output $d.Content.Text
